# The "Direction Info" sheet held a manually-curated classification table
# (Unit Number / Direction / Image Name) that duplicated what the
# GPT-based classifier now produces directly, so the sample data rows are
# removed, leaving just the empty header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data rows (2-10), shifting nothing below them (there is
# nothing below), so the sheet's used range shrinks to the header row.
$ws.Range("A2:C10").ClearContents()

# Clear the header labels themselves too, but keep their formatting/style.
$ws.Range("A1:C1").ClearContents()
